# Week 15 simulations added to the Receiving sheet (row 6 = D.Smythe, row 7 = J.Waddle)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Receiving")

# D.Smythe (row 6): Short Target, Short Comp, Deep Target, Deep Comp, RZ Target, RZ Comp
$ws.Range("C6").Value = 70
$ws.Range("D6").Value = 64
$ws.Range("E6").Value = 14
$ws.Range("F6").Value = 7
$ws.Range("G6").Value = 5
$ws.Range("H6").Value = 4

# J.Waddle (row 7): zeroed out
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0

# Make Receiving the active/selected sheet, with H7 selected (matches the saved view state)
$ws.Activate()
$ws.Range("H7").Select() | Out-Null
